# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list on Mon Jul 31 07:42:02 UTC 2023 with GitHub Actions".
#
# Column D ("Price") and E ("Volume(1h)") values are plain text in this sheet
# (prices use "." as a thousands separator, e.g. "29.420.30", and are not valid
# numbers). Any new value that WOULD parse as a plain number/float (e.g. "243.87",
# "1.000", "0.00000000127") is written with a leading single-quote so Excel keeps
# it as literal text instead of silently converting it to a Double (which would
# also strip meaningful trailing zeros like the ones in "1.000"/"9.000").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.420.30'
$ws.Range("E2").Value = '  +0.33%  '

# Row 3
$ws.Range("D3").Value = '1.868.70'
$ws.Range("E3").Value = '  -0.44%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").Value = '''243.87'
$ws.Range("E5").Value = '  +0.72%  '

# Row 6
$ws.Range("D6").Value = '''0.7072'
$ws.Range("E6").Value = '  -0.46%  '

# Row 7
$ws.Range("E7").Value = '  -0.07%  '

# Row 8
$ws.Range("D8").Value = '''0.07886'
$ws.Range("E8").Value = '  -1.36%  '

# Row 9
$ws.Range("E9").Value = '  -0.62%  '

# Row 10
$ws.Range("D10").Value = '''24.56'
$ws.Range("E10").Value = '  -1.60%  '

# Row 11
$ws.Range("D11").Value = '''0.07943'
$ws.Range("E11").Value = '  -4.00%  '

# Row 12
$ws.Range("D12").Value = '1.893.55'
$ws.Range("E12").Value = '  -0.08%  '

# Row 13
$ws.Range("D13").Value = '''5.211'
$ws.Range("E13").Value = '  -0.62%  '

# Row 14
$ws.Range("E14").Value = '  -1.04%  '

# Row 15
$ws.Range("D15").Value = '''0.7012'
$ws.Range("E15").Value = '  -1.23%  '

# Row 16
$ws.Range("D16").Value = '''6.507'
$ws.Range("E16").Value = '  +2.33%  '

# Row 17
$ws.Range("D17").Value = '29.473.88'
$ws.Range("E17").Value = '  +0.38%  '

# Row 18
$ws.Range("D18").Value = '''0.000008351'
$ws.Range("E18").Value = '  -1.87%  '

# Row 19
$ws.Range("D19").Value = '''252.27'
$ws.Range("E19").Value = '  +3.02%  '

# Row 20
$ws.Range("D20").Value = '2.131.26'
$ws.Range("E20").Value = '  -1.14%  '

# Row 21
$ws.Range("D21").Value = '''13.12'
$ws.Range("E21").Value = '  -1.03%  '

# Row 22
$ws.Range("D22").Value = '''1.000'
$ws.Range("E22").Value = '  -0.15%  '

# Row 23
$ws.Range("D23").Value = '''7.637'
$ws.Range("E23").Value = '  -1.77%  '

# Row 24
$ws.Range("E24").Value = '  -0.24%  '

# Row 25
$ws.Range("D25").Value = '''0.1554'
$ws.Range("E25").Value = '  -0.17%  '

# Row 26
$ws.Range("D26").Value = '''9.000'
$ws.Range("E26").Value = '  -0.64%  '

# Row 27
$ws.Range("D27").Value = '''161.60'
$ws.Range("E27").Value = '  -0.62%  '

# Row 28
$ws.Range("E28").Value = '  +1.01%  '

# Row 29
$ws.Range("D29").Value = '''1.501'
$ws.Range("E29").Value = '  -0.15%  '

# Row 30
$ws.Range("D30").Value = '''4.331'
$ws.Range("E30").Value = '  -1.87%  '

# Row 31
$ws.Range("D31").Value = '''4.252'
$ws.Range("E31").Value = '  -1.56%  '

# Row 32
$ws.Range("E32").Value = '  +1.48%  '

# Row 33
$ws.Range("D33").Value = '''0.05316'
$ws.Range("E33").Value = '  -1.00%  '

# Row 34
$ws.Range("D34").Value = '''1.895'
$ws.Range("E34").Value = '  -1.99%  '

# Row 35
$ws.Range("D35").Value = '''0.7487'
$ws.Range("E35").Value = '  -1.95%  '

# Row 36
$ws.Range("D36").Value = '''1.173'
$ws.Range("E36").Value = '  -0.83%  '

# Row 37
$ws.Range("D37").Value = '''2.714'
$ws.Range("E37").Value = '  +1.02%  '

# Row 38
$ws.Range("D38").Value = '''0.01887'
$ws.Range("E38").Value = '  +0.37%  '

# Row 39
$ws.Range("D39").Value = '1.274.57'
$ws.Range("E39").Value = '  +1.35%  '

# Row 40
$ws.Range("D40").Value = '''2.770'
$ws.Range("E40").Value = '  +0.69%  '

# Row 41
$ws.Range("D41").Value = '''0.8930'
$ws.Range("E41").Value = '  -1.79%  '

# Row 42
$ws.Range("E42").Value = '  -6.82%  '

# Row 43
$ws.Range("E43").Value = '  -3.20%  '

# Row 44
$ws.Range("D44").Value = '''71.35'
$ws.Range("E44").Value = '  -3.79%  '

# Row 45
$ws.Range("E45").Value = '  -0.11%  '

# Row 46
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '''0.00000000127'
$ws.Range("E46").Value = '  -3.65%  '

# Row 47
$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").Value = '2.030.49'
$ws.Range("E47").Value = '  -1.19%  '

# Row 48
$ws.Range("E48").Value = '  -0.02%  '

# Row 49
$ws.Range("D49").Value = '''9.555'
$ws.Range("E49").Value = '  +1.28%  '

# Row 50
$ws.Range("D50").Value = '''0.5182'
$ws.Range("E50").Value = '  -0.83%  '

# Row 51
$ws.Range("D51").Value = '''0.4306'
$ws.Range("E51").Value = '  -1.41%  '
